$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.513.58"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.563.74"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'211.64"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'46.39"
$ws.Range("E8").Value = "  +5.45%  "
$ws.Range("D9").Value = "'24.10"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "1.787.10"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "1.544.85"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").Value = "28.508.65"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "'3.68"
$ws.Range("E17").Value = "  -3.18%  "
$ws.Range("D18").Value = "'62.16"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "'228.99"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").Value = "'7.34"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'3.87"
$ws.Range("E23").Value = "  -6.19%  "
$ws.Range("D24").Value = "'9.14"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("D25").Value = "'2.09"
$ws.Range("E25").Value = "  +6.59%  "
$ws.Range("D26").Value = "'150.39"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "'14.98"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").Value = "'6.44"
$ws.Range("E29").Value = "  -4.15%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "'0.0466"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").Value = "'1.10"
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("D33").Value = "'3.20"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").Value = "1.396.38"
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("D38").Value = "'2.36"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "'2.58"
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "'0.536"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'1.89"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("E44").Value = "  -4.04%  "
$ws.Range("E45").Value = "  -4.44%  "
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").Value = "'62.68"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").Value = "1.700.07"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("D49").Value = "'86.11"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0103"
$ws.Range("E50").Value = "  -4.02%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0₆0103"
$ws.Range("E51").Value = "  -0.08%  "
